# The document contains three <id>...</id> markers, each currently split
# across three separate runs: "<id>", the old id value ("p004v_aN"), and
# "</id>". Each needs to become a single run whose text reads
# "<id>p004v_N</id>" (dropping the "a" from the old id and renumbering).
#
# Using Find & Replace with cross-run matching both updates the text and
# collapses the three runs into one, adopting the formatting (Courier New,
# color 7f6000, sz 18) of the first of the three runs - matching the target
# edit exactly.

$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p004v_a1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p004v_1</id>", 2)
$d.Content.Find.Execute("<id>p004v_a2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p004v_2</id>", 2)
$d.Content.Find.Execute("<id>p004v_a3</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p004v_3</id>", 2)
